$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-4 (columns A:T)
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics
$rows = @(
    @("ECs",  "Cck", "Cckbr", "FAPs", 1, 0.3333333333333333, 0.1794766666666666, 0.53843,        0.0361256439470005,  0.03612564394700051, 3, 1, 0.83582, 2.50746, 1, 1, 0.1500101875333333, 1.3500916878,     0.0361256439470005,  0.03612564394700051),
    @("FAPs", "Cck", "Cckbr", "FAPs", 1, 0.3333333333333333, 0.3771426666666667, 1.131428,       0.07591249573698883, 0.07591249573698884, 3, 1, 0.83582, 2.50746, 1, 1, 0.3152233836533334, 2.83701045288,    0.07591249573698883, 0.07591249573698884),
    @("sCs",  "Cck", "Cckbr", "FAPs", 3, 1,                  4.411504333333334, 13.234513,       0.8879618603160108,  0.8879618603160107,  3, 1, 0.83582, 2.50746, 1, 1, 3.687223551886667, 33.18501196698001, 0.8879618603160108,  0.8879618603160107)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
